$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.333
$ws.Range("B4").Value = 0.0155
$ws.Range("B5").Value = 0.0691
$ws.Range("B6").Value = 0.1664
$ws.Range("B7").Value = -0.1477
$ws.Range("B8").Value = -0.2364
$ws.Range("B9").Value = 0.003
$ws.Range("B10").Value = 0.0167
$ws.Range("B12").Value = 0.5516
$ws.Range("B13").Value = 0.0097
$ws.Range("B14").Value = 0.3531
$ws.Range("B15").Value = 0.073
$ws.Range("B16").Value = -0.0207
$ws.Range("B17").Value = 0.0058
$ws.Range("B18").Value = -0.0344
$ws.Range("B20").Value = -0.0002
$ws.Range("B21").Value = -0.1057
$ws.Range("B22").Value = 0.0003
$ws.Range("B23").Value = 0.2896
$ws.Range("B24").Value = 0.0372
